# Update the "想去人数" (interest count, column F) values on each sheet
# to match the latest scrape, per commit "Update gh-pages to output generated at 456a3b4".
$wb = $excel.ActiveWorkbook

# Sheet 1: 展览
$ws = $wb.Worksheets.Item(1)
$ws.Range("F2").Value = 1549  # was 1548
$ws.Range("F5").Value = 8426  # was 8378
$ws.Range("F6").Value = 228  # was 227
$ws.Range("F7").Value = 90  # was 89
$ws.Range("F8").Value = 1223  # was 1219
$ws.Range("F10").Value = 227  # was 223
$ws.Range("F13").Value = 82  # was 80
$ws.Range("F14").Value = 265  # was 264
$ws.Range("F17").Value = 1364  # was 1362
$ws.Range("F21").Value = 1312  # was 1311
$ws.Range("F25").Value = 54  # was 53
$ws.Range("F27").Value = 255  # was 253
$ws.Range("F28").Value = 1052  # was 1051
$ws.Range("F30").Value = 22  # was 21
$ws.Range("F32").Value = 159  # was 157
$ws.Range("F42").Value = 655  # was 654

# Sheet 2: 演出
$ws = $wb.Worksheets.Item(2)
$ws.Range("F12").Value = 216  # was 217
$ws.Range("F21").Value = 50  # was 49
$ws.Range("F26").Value = 1019  # was 1018
$ws.Range("F27").Value = 43  # was 40
$ws.Range("F28").Value = 629  # was 628

# Sheet 3: 本地生活
$ws = $wb.Worksheets.Item(3)
$ws.Range("F7").Value = 248  # was 247
$ws.Range("F8").Value = 119  # was 118
$ws.Range("F9").Value = 1901  # was 1893
$ws.Range("F10").Value = 2878  # was 2870

# Sheet 4: 全部类型
$ws = $wb.Worksheets.Item(4)
$ws.Range("F3").Value = 1549  # was 1548
$ws.Range("F7").Value = 8426  # was 8378
$ws.Range("F8").Value = 248  # was 247
$ws.Range("F9").Value = 119  # was 118
$ws.Range("F10").Value = 228  # was 227
$ws.Range("F11").Value = 1901  # was 1893
$ws.Range("F12").Value = 2878  # was 2870
$ws.Range("F13").Value = 216  # was 217
$ws.Range("F14").Value = 90  # was 89
$ws.Range("F15").Value = 1223  # was 1219
$ws.Range("F17").Value = 227  # was 223
$ws.Range("F20").Value = 265  # was 264
$ws.Range("F21").Value = 1364  # was 1362
$ws.Range("F25").Value = 1312  # was 1311
$ws.Range("F27").Value = 54  # was 53
$ws.Range("F31").Value = 50  # was 49
$ws.Range("F32").Value = 22  # was 21
$ws.Range("F35").Value = 159  # was 157
$ws.Range("F36").Value = 1019  # was 1018
$ws.Range("F39").Value = 629  # was 628
$ws.Range("F43").Value = 655  # was 654
